$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Relations sheet: correct the targets of relations r and s so the
#    typing lines up with the new Rule (r: A->C, s: A->B, was A->B/A->C)
# ---------------------------------------------------------------------
$wsRelations = $wb.Worksheets.Item("Relations")
$wsRelations.Range("D3").Value = "C"
$wsRelations.Range("D4").Value = "B"
$wsRelations.Range("D4").Select()

# ---------------------------------------------------------------------
# 2. Rules sheet: move the [Rule] header into A1 and drop the stray
#    E6 cell, leaving a clean A1:C3 table.
# ---------------------------------------------------------------------
$wsRules = $wb.Worksheets.Item("Rules")
$wsRules.Range("E6").ClearContents()
$wsRules.Range("A1").Value = "[Rule]"
$wsRules.Range("B3:C3").Select()

# ---------------------------------------------------------------------
# 3. Terms sheet: rename t4's representation (shared string order must
#    match: "r ISC s;t" is introduced here, before "=" and before the
#    combined equation string below).
# ---------------------------------------------------------------------
$wsTerms = $wb.Worksheets.Item("Terms")
$wsTerms.Range("B6").Value = "r ISC s;t"

# ---------------------------------------------------------------------
# 4. Binary Terms sheet: add the binary term representing the rule's
#    equation: r1 = (t5 = t4), i.e. operator "=" between lhs t5, rhs t4.
# ---------------------------------------------------------------------
$wsBinaryTerms = $wb.Worksheets.Item("Binary Terms")
$wsBinaryTerms.Range("A8").Value = "r1"
$wsBinaryTerms.Range("C8").NumberFormat = "@"
$wsBinaryTerms.Range("C8").Value = "t5"
$wsBinaryTerms.Range("D8").NumberFormat = "@"
$wsBinaryTerms.Range("D8").Value = "t4"
$wsBinaryTerms.Range("B8").Value = "'="
$wsBinaryTerms.Range("B8").NumberFormat = "@"
$wsBinaryTerms.Range("B9").Select()

# ---------------------------------------------------------------------
# 5. Terms sheet: add the new term r1 whose representation is the
#    rule's equation (introduces the final new shared string).
# ---------------------------------------------------------------------
$wsTerms.Range("A8").Value = "r1"
$wsTerms.Range("B8").Value = "s;t = r ISC s;t"

# ---------------------------------------------------------------------
# 6. Restore the Terms sheet as the active tab/selection, matching the
#    original workbook (activeTab=4, Terms tabSelected).
# ---------------------------------------------------------------------
$wsTerms.Range("B9").Select()
